$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for Neutrophils as sending cluster (rows 11-13); dataset updated with new TPM values
$ws.Rows("11:13").Delete()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.1450045
$ws.Range("H2").Value = 14.290009
$ws.Range("I2").Value = 0.8119737125238713
$ws.Range("J2").Value = 0.7990590344890214
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.0421005
$ws.Range("N2").Value = 0.084201
$ws.Range("O2").Value = 0.02355433709362141
$ws.Range("P2").Value = 0.02354456735134313
$ws.Range("Q2").Value = 0.30080826195225
$ws.Range("R2").Value = 1.203233047809
$ws.Range("S2").Value = 0.01912550253594651
$ws.Range("T2").Value = 0.01881349925522597

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.1450045
$ws.Range("H3").Value = 14.290009
$ws.Range("I3").Value = 0.8119737125238713
$ws.Range("J3").Value = 0.7990590344890214
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.001483333333333333
$ws.Range("N3").Value = 0.00445
$ws.Range("O3").Value = 0.000829893548901757
$ws.Range("P3").Value = 0.001244323995124487
$ws.Range("Q3").Value = 0.01059842334166667
$ws.Range("R3").Value = 0.06359054004999999
$ws.Range("S3").Value = 0.0006738517459013706
$ws.Range("T3").Value = 0.0009942883301356942

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.1450045
$ws.Range("H4").Value = 14.290009
$ws.Range("I4").Value = 0.8119737125238713
$ws.Range("J4").Value = 0.7990590344890214
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.743794
$ws.Range("N4").Value = 3.487588
$ws.Range("O4").Value = 0.9756157693574768
$ws.Range("P4").Value = 0.9752111086535323
$ws.Range("Q4").Value = 12.459415977073
$ws.Range("R4").Value = 49.837663908292
$ws.Range("S4").Value = 0.7921743582420234
$ws.Range("T4").Value = 0.7792512469036597

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ntf3"
$ws.Range("C5").Value = "Ntrk3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.284443
$ws.Range("H5").Value = 0.853329
$ws.Range("I5").Value = 0.0323247156403369
$ws.Range("J5").Value = 0.04771587245616726
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.0421005
$ws.Range("N5").Value = 0.084201
$ws.Range("O5").Value = 0.02355433709362141
$ws.Range("P5").Value = 0.02354456735134313
$ws.Range("Q5").Value = 0.0119751925215
$ws.Range("R5").Value = 0.071851155129
$ws.Range("S5").Value = 0.0007613872486479515
$ws.Range("T5").Value = 0.001123449572772328

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntf3"
$ws.Range("C6").Value = "Ntrk3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.284443
$ws.Range("H6").Value = 0.853329
$ws.Range("I6").Value = 0.0323247156403369
$ws.Range("J6").Value = 0.04771587245616726
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.001483333333333333
$ws.Range("N6").Value = 0.00445
$ws.Range("O6").Value = 0.000829893548901757
$ws.Range("P6").Value = 0.001244323995124487
$ws.Range("Q6").Value = 0.0004219237833333333
$ws.Range("R6").Value = 0.00379731405
$ws.Range("S6").Value = [double]"2.682607297999932e-05"
$ws.Range("T6").Value = [double]"5.93740050455085e-05"

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntf3"
$ws.Range("C7").Value = "Ntrk3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.284443
$ws.Range("H7").Value = 0.853329
$ws.Range("I7").Value = 0.0323247156403369
$ws.Range("J7").Value = 0.04771587245616726
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.743794
$ws.Range("N7").Value = 3.487588
$ws.Range("O7").Value = 0.9756157693574768
$ws.Range("P7").Value = 0.9752111086535323
$ws.Range("Q7").Value = 0.496009996742
$ws.Range("R7").Value = 2.976059980452
$ws.Range("S7").Value = 0.03153650231870895
$ws.Range("T7").Value = 0.04653304887834941

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ntf3"
$ws.Range("C8").Value = "Ntrk3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.370104
$ws.Range("H8").Value = 2.740208
$ws.Range("I8").Value = 0.1557015718357919
$ws.Range("J8").Value = 0.1532250930548114
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.0421005
$ws.Range("N8").Value = 0.084201
$ws.Range("O8").Value = 0.02355433709362141
$ws.Range("P8").Value = 0.02354456735134313
$ws.Range("Q8").Value = 0.057682063452
$ws.Range("R8").Value = 0.230728253808
$ws.Range("S8").Value = 0.003667447309026951
$ws.Range("T8").Value = 0.003607618523344824

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ntf3"
$ws.Range("C9").Value = "Ntrk3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.370104
$ws.Range("H9").Value = 2.740208
$ws.Range("I9").Value = 0.1557015718357919
$ws.Range("J9").Value = 0.1532250930548114
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.001483333333333333
$ws.Range("N9").Value = 0.00445
$ws.Range("O9").Value = 0.000829893548901757
$ws.Range("P9").Value = 0.001244323995124487
$ws.Range("Q9").Value = 0.002032320933333333
$ws.Range("R9").Value = 0.0121939256
$ws.Range("S9").Value = 0.0001292157300203872
$ws.Range("T9").Value = 0.0001906616599432842

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntf3"
$ws.Range("C10").Value = "Ntrk3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.370104
$ws.Range("H10").Value = 2.740208
$ws.Range("I10").Value = 0.1557015718357919
$ws.Range("J10").Value = 0.1532250930548114
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.743794
$ws.Range("N10").Value = 3.487588
$ws.Range("O10").Value = 0.9756157693574768
$ws.Range("P10").Value = 0.9752111086535323
$ws.Range("Q10").Value = 2.389179134576
$ws.Range("R10").Value = 9.556716538304
$ws.Range("S10").Value = 0.1519049087967445
$ws.Range("T10").Value = 0.1494268128715233
